$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 18
$ws.Range("B18").Value = "[-, 'ELM-2NA-Tecnologia da Soldagem', -, -]"
$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "[-, -, 'MEC-2NB-Soldagem', -]"
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = "-"

# Row 19
$ws.Range("B19").Value = "[-, 'ELM-2NA-Tecnologia da Soldagem', -, -]"
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "[-, -, 'MEC-2NB-Soldagem', -]"
$ws.Range("E19").Value = "-"
$ws.Range("F19").Value = "-"

# Row 20
$ws.Range("B20").Value = "[-, 'ELM-2NA-Tecnologia da Soldagem', -, -]"
$ws.Range("D20").Value = "[-, -, 'MEC-2NB-Soldagem', -]"
$ws.Range("F20").Value = "-"

# Row 21
$ws.Range("B21").Value = "[-, 'ELM-2NA-Tecnologia da Soldagem', -, -]"
$ws.Range("D21").Value = "[-, -, 'MEC-2NB-Soldagem', -]"
$ws.Range("F21").Value = "-"
